# Applies the diff for 9week/2019775054.pptx:
#  - Slide 1, TextBox "게임프로그래밍": widen/re-center the title box
#    (off x 4188228 -> 3564294, ext cx 4636774 -> 5260708; right edge
#    stays fixed at 8825002 EMU).
#  - Slide 10, TextBox "1. 게임 실행": widen the box (ext cx 2770310 ->
#    3029997) and renumber "1." -> "2.".
#  - Slide 11, TextBox "4. 출처": renumber "4." -> "3.".
#
# NOTE on the literal Points values below: Shape.Left/Top/Width/Height
# are COM `Single` (32-bit float) properties. 1 pt = 12700 EMU, so an
# EMU target is first divided by 12700 and then narrowed to the nearest
# float32 representable value that round-trips (via truncation) back to
# the exact target EMU, instead of the naive `emu/12700`, which can be
# off by 1 EMU after the float32 round-trip.

$p = $ppt.ActivePresentation

# --- Slide 1: title textbox "게임프로그래밍" -------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$sh1.Left = 280.653076171875
$sh1.Width = 414.22900390625

# --- Slide 10: "1. 게임 실행" -> "2. 게임 실행" ------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(3)
$sh10.Width = 238.5824432373047
$tr10 = $sh10.TextFrame.TextRange
$num10 = $tr10.Characters(1, 2)
$num10.Text = "2."

# --- Slide 11: "4. 출처" -> "3. 출처" ---------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(3)
$tr11 = $sh11.TextFrame.TextRange
$num11 = $tr11.Characters(1, 3)
$num11.Text = "3. "
